# Applies the "Updated cryptos list" GitHub Actions commit.
#
# The sheet stores Price (D) / Volume(1h) (E) figures, and a couple of
# Coin/Link/Price rows (42/43 and 45/46), as literal text (inline strings),
# not numbers -- many "prices" use dotted thousands separators
# (e.g. "44.059.57") that are not valid numeric literals, and the percent
# cells keep two leading/trailing spaces of padding. To stop Excel from
# "helpfully" reinterpreting plain numeric-looking text (like "1.00" or
# "0.100") as a number and dropping the significant trailing zero, each
# write temporarily forces the cell to Text format, then restores the
# original "Normal" style so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "43.963.44"
Set-TextValue "E2" "  +0.03%  "
Set-TextValue "D3" "2.352.48"
Set-TextValue "E3" "  -0.43%  "
Set-TextValue "E4" "  +0.48%  "
Set-TextValue "D5" "0.689"
Set-TextValue "E5" "  +4.29%  "
Set-TextValue "D6" "239.99"
Set-TextValue "E6" "  +2.09%  "
Set-TextValue "D7" "75.82"
Set-TextValue "E7" "  +3.00%  "
Set-TextValue "E8" "  +0.05%  "
Set-TextValue "D9" "0.612"
Set-TextValue "E9" "  +15.16%  "
Set-TextValue "D10" "0.100"
Set-TextValue "E10" "  +1.78%  "
Set-TextValue "D11" "57.28"
Set-TextValue "E11" "  +0.65%  "
Set-TextValue "D12" "33.05"
Set-TextValue "E12" "  +16.23%  "
Set-TextValue "D13" "7.39"
Set-TextValue "E13" "  +11.48%  "
Set-TextValue "D14" "0.108"
Set-TextValue "E14" "  +1.73%  "
Set-TextValue "D15" "2.703.20"
Set-TextValue "E15" "  -0.33%  "
Set-TextValue "D16" "16.59"
Set-TextValue "E16" "  -1.25%  "
Set-TextValue "D17" "0.918"
Set-TextValue "E17" "  +3.72%  "
Set-TextValue "D18" "2.352.25"
Set-TextValue "E18" "  -0.17%  "
Set-TextValue "D19" "43.896.16"
Set-TextValue "E19" "  +0.16%  "
Set-TextValue "E20" "  +0.84%  "
Set-TextValue "D21" "6.66"
Set-TextValue "E21" "  +5.16%  "
Set-TextValue "D22" "77.32"
Set-TextValue "E22" "  +1.71%  "
Set-TextValue "D23" "257.90"
Set-TextValue "E23" "  +2.60%  "
Set-TextValue "D24" "0.999"
Set-TextValue "E24" "  -0.08%  "
Set-TextValue "E25" "  -1.25%  "
Set-TextValue "D26" "2.51"
Set-TextValue "E26" "  +1.02%  "
Set-TextValue "E27" "  +15.83%  "
Set-TextValue "D28" "10.82"
Set-TextValue "E28" "  +5.32%  "
Set-TextValue "D29" "2.28"
Set-TextValue "E29" "  +1.15%  "
Set-TextValue "D30" "22.94"
Set-TextValue "E30" "  +1.92%  "
Set-TextValue "D31" "174.42"
Set-TextValue "E31" "  +0.80%  "
Set-TextValue "E32" "  -4.04%  "
Set-TextValue "E33" "  +3.85%  "
Set-TextValue "D34" "5.34"
Set-TextValue "E34" "  +3.43%  "
Set-TextValue "D35" "0.0758"
Set-TextValue "E35" "  +7.04%  "
Set-TextValue "D36" "5.38"
Set-TextValue "E36" "  +4.98%  "
Set-TextValue "D37" "3.74"
Set-TextValue "E37" "  -0.70%  "
Set-TextValue "D38" "2.40"
Set-TextValue "E38" "  -1.71%  "
Set-TextValue "D39" "6.39"
Set-TextValue "E39" "  -0.95%  "
Set-TextValue "D40" "0.0281"
Set-TextValue "E40" "  +6.43%  "
Set-TextValue "E41" "  +21.07%  "
Set-TextValue "B42" "Cronos"
Set-TextValue "C42" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D42" "0.107"
Set-TextValue "E42" "  +11.18%  "
Set-TextValue "B43" "FraxShare"
Set-TextValue "C43" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D43" "9.11"
Set-TextValue "E43" "  +2.32%  "
Set-TextValue "D44" "19.05"
Set-TextValue "E44" "  -1.71%  "
Set-TextValue "B45" "BinanceUSD"
Set-TextValue "C45" "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue "D45" "1.00"
Set-TextValue "E45" "  +0.17%  "
Set-TextValue "B46" "FTXToken"
Set-TextValue "C46" "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue "D46" "4.83"
Set-TextValue "E46" "  +8.82%  "
Set-TextValue "D47" "2.52"
Set-TextValue "E47" "  +9.53%  "
Set-TextValue "D48" "1.25"
Set-TextValue "E48" "  +2.43%  "
Set-TextValue "D49" "101.40"
Set-TextValue "E49" "  +2.39%  "
Set-TextValue "D50" "1.18"
Set-TextValue "E50" "  +0.40%  "
Set-TextValue "D51" "55.68"
Set-TextValue "E51" "  +7.39%  "
